$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("PIR")
$ws.Range("A227:A240").NumberFormat = "@"

$ws.Cells.Item(227, 1).Value = "2026-01-28"
$ws.Cells.Item(227, 2).Value = "16:29:24"
$ws.Cells.Item(227, 3).Value = "16:00"
$ws.Cells.Item(227, 4).Value = "Bathroom"
$ws.Cells.Item(227, 5).Value = "No Motion"
$ws.Cells.Item(227, 6).Value = "Inactive"

$ws.Cells.Item(228, 1).Value = "2026-01-28"
$ws.Cells.Item(228, 2).Value = "16:29:24"
$ws.Cells.Item(228, 3).Value = "16:00"
$ws.Cells.Item(228, 4).Value = "Bathroom"
$ws.Cells.Item(228, 5).Value = "No Motion"
$ws.Cells.Item(228, 6).Value = "Inactive"

$ws.Cells.Item(229, 1).Value = "2026-01-28"
$ws.Cells.Item(229, 2).Value = "16:29:29"
$ws.Cells.Item(229, 3).Value = "16:00"
$ws.Cells.Item(229, 4).Value = "Bathroom"
$ws.Cells.Item(229, 5).Value = "No Motion"
$ws.Cells.Item(229, 6).Value = "Inactive"

$ws.Cells.Item(230, 1).Value = "2026-01-28"
$ws.Cells.Item(230, 2).Value = "16:29:34"
$ws.Cells.Item(230, 3).Value = "16:00"
$ws.Cells.Item(230, 4).Value = "Bathroom"
$ws.Cells.Item(230, 5).Value = "No Motion"
$ws.Cells.Item(230, 6).Value = "Inactive"

$ws.Cells.Item(231, 1).Value = "2026-01-28"
$ws.Cells.Item(231, 2).Value = "16:29:39"
$ws.Cells.Item(231, 3).Value = "16:00"
$ws.Cells.Item(231, 4).Value = "Bathroom"
$ws.Cells.Item(231, 5).Value = "No Motion"
$ws.Cells.Item(231, 6).Value = "Inactive"

$ws.Cells.Item(232, 1).Value = "2026-01-28"
$ws.Cells.Item(232, 2).Value = "16:29:45"
$ws.Cells.Item(232, 3).Value = "16:00"
$ws.Cells.Item(232, 4).Value = "Bathroom"
$ws.Cells.Item(232, 5).Value = "No Motion"
$ws.Cells.Item(232, 6).Value = "Inactive"

$ws.Cells.Item(233, 1).Value = "2026-01-28"
$ws.Cells.Item(233, 2).Value = "16:29:49"
$ws.Cells.Item(233, 3).Value = "16:00"
$ws.Cells.Item(233, 4).Value = "Bathroom"
$ws.Cells.Item(233, 5).Value = "No Motion"
$ws.Cells.Item(233, 6).Value = "Inactive"

$ws.Cells.Item(234, 1).Value = "2026-01-28"
$ws.Cells.Item(234, 2).Value = "16:29:54"
$ws.Cells.Item(234, 3).Value = "16:00"
$ws.Cells.Item(234, 4).Value = "Bathroom"
$ws.Cells.Item(234, 5).Value = "No Motion"
$ws.Cells.Item(234, 6).Value = "Inactive"

$ws.Cells.Item(235, 1).Value = "2026-01-28"
$ws.Cells.Item(235, 2).Value = "16:29:59"
$ws.Cells.Item(235, 3).Value = "16:00"
$ws.Cells.Item(235, 4).Value = "Bathroom"
$ws.Cells.Item(235, 5).Value = "No Motion"
$ws.Cells.Item(235, 6).Value = "Inactive"

$ws.Cells.Item(236, 1).Value = "2026-01-28"
$ws.Cells.Item(236, 2).Value = "16:30:05"
$ws.Cells.Item(236, 3).Value = "16:00"
$ws.Cells.Item(236, 4).Value = "Bathroom"
$ws.Cells.Item(236, 5).Value = "No Motion"
$ws.Cells.Item(236, 6).Value = "Inactive"

$ws.Cells.Item(237, 1).Value = "2026-01-28"
$ws.Cells.Item(237, 2).Value = "16:30:09"
$ws.Cells.Item(237, 3).Value = "16:00"
$ws.Cells.Item(237, 4).Value = "Bathroom"
$ws.Cells.Item(237, 5).Value = "No Motion"
$ws.Cells.Item(237, 6).Value = "Inactive"

$ws.Cells.Item(238, 1).Value = "2026-01-28"
$ws.Cells.Item(238, 2).Value = "16:30:14"
$ws.Cells.Item(238, 3).Value = "16:00"
$ws.Cells.Item(238, 4).Value = "Bathroom"
$ws.Cells.Item(238, 5).Value = "No Motion"
$ws.Cells.Item(238, 6).Value = "Inactive"

$ws.Cells.Item(239, 1).Value = "2026-01-28"
$ws.Cells.Item(239, 2).Value = "16:30:19"
$ws.Cells.Item(239, 3).Value = "16:00"
$ws.Cells.Item(239, 4).Value = "Bathroom"
$ws.Cells.Item(239, 5).Value = "No Motion"
$ws.Cells.Item(239, 6).Value = "Inactive"

$ws.Cells.Item(240, 1).Value = "2026-01-28"
$ws.Cells.Item(240, 2).Value = "16:30:25"
$ws.Cells.Item(240, 3).Value = "16:00"
$ws.Cells.Item(240, 4).Value = "Bathroom"
$ws.Cells.Item(240, 5).Value = "No Motion"
$ws.Cells.Item(240, 6).Value = "Inactive"


$ws = $wb.Worksheets.Item("Humidity")
$ws.Range("A224:A234").NumberFormat = "@"
$ws.Range("E224:E234").NumberFormat = "@"

$ws.Cells.Item(224, 1).Value = "2026-01-28"
$ws.Cells.Item(224, 2).Value = "16:29:27"
$ws.Cells.Item(224, 3).Value = "16:00"
$ws.Cells.Item(224, 4).Value = "Bathroom"
$ws.Cells.Item(224, 5).Value = "87.8%"
$ws.Cells.Item(224, 6).Value = "Active"

$ws.Cells.Item(225, 1).Value = "2026-01-28"
$ws.Cells.Item(225, 2).Value = "16:29:31"
$ws.Cells.Item(225, 3).Value = "16:00"
$ws.Cells.Item(225, 4).Value = "Bathroom"
$ws.Cells.Item(225, 5).Value = "86.9%"
$ws.Cells.Item(225, 6).Value = "Active"

$ws.Cells.Item(226, 1).Value = "2026-01-28"
$ws.Cells.Item(226, 2).Value = "16:29:43"
$ws.Cells.Item(226, 3).Value = "16:00"
$ws.Cells.Item(226, 4).Value = "Bathroom"
$ws.Cells.Item(226, 5).Value = "87.9%"
$ws.Cells.Item(226, 6).Value = "Active"

$ws.Cells.Item(227, 1).Value = "2026-01-28"
$ws.Cells.Item(227, 2).Value = "16:29:51"
$ws.Cells.Item(227, 3).Value = "16:00"
$ws.Cells.Item(227, 4).Value = "Bathroom"
$ws.Cells.Item(227, 5).Value = "87.0%"
$ws.Cells.Item(227, 6).Value = "Active"

$ws.Cells.Item(228, 1).Value = "2026-01-28"
$ws.Cells.Item(228, 2).Value = "16:29:55"
$ws.Cells.Item(228, 3).Value = "16:00"
$ws.Cells.Item(228, 4).Value = "Bathroom"
$ws.Cells.Item(228, 5).Value = "87.9%"
$ws.Cells.Item(228, 6).Value = "Active"

$ws.Cells.Item(229, 1).Value = "2026-01-28"
$ws.Cells.Item(229, 2).Value = "16:30:00"
$ws.Cells.Item(229, 3).Value = "16:00"
$ws.Cells.Item(229, 4).Value = "Bathroom"
$ws.Cells.Item(229, 5).Value = "87.9%"
$ws.Cells.Item(229, 6).Value = "Active"

$ws.Cells.Item(230, 1).Value = "2026-01-28"
$ws.Cells.Item(230, 2).Value = "16:30:03"
$ws.Cells.Item(230, 3).Value = "16:00"
$ws.Cells.Item(230, 4).Value = "Bathroom"
$ws.Cells.Item(230, 5).Value = "88.0%"
$ws.Cells.Item(230, 6).Value = "Active"

$ws.Cells.Item(231, 1).Value = "2026-01-28"
$ws.Cells.Item(231, 2).Value = "16:30:08"
$ws.Cells.Item(231, 3).Value = "16:00"
$ws.Cells.Item(231, 4).Value = "Bathroom"
$ws.Cells.Item(231, 5).Value = "87.9%"
$ws.Cells.Item(231, 6).Value = "Active"

$ws.Cells.Item(232, 1).Value = "2026-01-28"
$ws.Cells.Item(232, 2).Value = "16:30:16"
$ws.Cells.Item(232, 3).Value = "16:00"
$ws.Cells.Item(232, 4).Value = "Bathroom"
$ws.Cells.Item(232, 5).Value = "87.9%"
$ws.Cells.Item(232, 6).Value = "Active"

$ws.Cells.Item(233, 1).Value = "2026-01-28"
$ws.Cells.Item(233, 2).Value = "16:30:20"
$ws.Cells.Item(233, 3).Value = "16:00"
$ws.Cells.Item(233, 4).Value = "Bathroom"
$ws.Cells.Item(233, 5).Value = "87.9%"
$ws.Cells.Item(233, 6).Value = "Active"

$ws.Cells.Item(234, 1).Value = "2026-01-28"
$ws.Cells.Item(234, 2).Value = "16:30:24"
$ws.Cells.Item(234, 3).Value = "16:00"
$ws.Cells.Item(234, 4).Value = "Bathroom"
$ws.Cells.Item(234, 5).Value = "87.0%"
$ws.Cells.Item(234, 6).Value = "Active"


$ws = $wb.Worksheets.Item("Temperature")
$ws.Range("A224:A234").NumberFormat = "@"

$ws.Cells.Item(224, 1).Value = "2026-01-28"
$ws.Cells.Item(224, 2).Value = "16:29:28"
$ws.Cells.Item(224, 3).Value = "16:00"
$ws.Cells.Item(224, 4).Value = "Bathroom"
$ws.Cells.Item(224, 5).Value = "22.8C"
$ws.Cells.Item(224, 6).Value = "Active"

$ws.Cells.Item(225, 1).Value = "2026-01-28"
$ws.Cells.Item(225, 2).Value = "16:29:32"
$ws.Cells.Item(225, 3).Value = "16:00"
$ws.Cells.Item(225, 4).Value = "Bathroom"
$ws.Cells.Item(225, 5).Value = "22.8C"
$ws.Cells.Item(225, 6).Value = "Active"

$ws.Cells.Item(226, 1).Value = "2026-01-28"
$ws.Cells.Item(226, 2).Value = "16:29:44"
$ws.Cells.Item(226, 3).Value = "16:00"
$ws.Cells.Item(226, 4).Value = "Bathroom"
$ws.Cells.Item(226, 5).Value = "22.8C"
$ws.Cells.Item(226, 6).Value = "Active"

$ws.Cells.Item(227, 1).Value = "2026-01-28"
$ws.Cells.Item(227, 2).Value = "16:29:52"
$ws.Cells.Item(227, 3).Value = "16:00"
$ws.Cells.Item(227, 4).Value = "Bathroom"
$ws.Cells.Item(227, 5).Value = "22.8C"
$ws.Cells.Item(227, 6).Value = "Active"

$ws.Cells.Item(228, 1).Value = "2026-01-28"
$ws.Cells.Item(228, 2).Value = "16:29:56"
$ws.Cells.Item(228, 3).Value = "16:00"
$ws.Cells.Item(228, 4).Value = "Bathroom"
$ws.Cells.Item(228, 5).Value = "22.8C"
$ws.Cells.Item(228, 6).Value = "Active"

$ws.Cells.Item(229, 1).Value = "2026-01-28"
$ws.Cells.Item(229, 2).Value = "16:30:00"
$ws.Cells.Item(229, 3).Value = "16:00"
$ws.Cells.Item(229, 4).Value = "Bathroom"
$ws.Cells.Item(229, 5).Value = "22.8C"
$ws.Cells.Item(229, 6).Value = "Active"

$ws.Cells.Item(230, 1).Value = "2026-01-28"
$ws.Cells.Item(230, 2).Value = "16:30:04"
$ws.Cells.Item(230, 3).Value = "16:00"
$ws.Cells.Item(230, 4).Value = "Bathroom"
$ws.Cells.Item(230, 5).Value = "22.8C"
$ws.Cells.Item(230, 6).Value = "Active"

$ws.Cells.Item(231, 1).Value = "2026-01-28"
$ws.Cells.Item(231, 2).Value = "16:30:08"
$ws.Cells.Item(231, 3).Value = "16:00"
$ws.Cells.Item(231, 4).Value = "Bathroom"
$ws.Cells.Item(231, 5).Value = "22.8C"
$ws.Cells.Item(231, 6).Value = "Active"

$ws.Cells.Item(232, 1).Value = "2026-01-28"
$ws.Cells.Item(232, 2).Value = "16:30:16"
$ws.Cells.Item(232, 3).Value = "16:00"
$ws.Cells.Item(232, 4).Value = "Bathroom"
$ws.Cells.Item(232, 5).Value = "22.8C"
$ws.Cells.Item(232, 6).Value = "Active"

$ws.Cells.Item(233, 1).Value = "2026-01-28"
$ws.Cells.Item(233, 2).Value = "16:30:21"
$ws.Cells.Item(233, 3).Value = "16:00"
$ws.Cells.Item(233, 4).Value = "Bathroom"
$ws.Cells.Item(233, 5).Value = "22.8C"
$ws.Cells.Item(233, 6).Value = "Active"

$ws.Cells.Item(234, 1).Value = "2026-01-28"
$ws.Cells.Item(234, 2).Value = "16:30:24"
$ws.Cells.Item(234, 3).Value = "16:00"
$ws.Cells.Item(234, 4).Value = "Bathroom"
$ws.Cells.Item(234, 5).Value = "22.8C"
$ws.Cells.Item(234, 6).Value = "Active"

